$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix the "two parallel starting operations" issue: the first operation of
#    a job must be marked "S" (start), not "P1" (a second parallel start).
#    O13 (row4), O15 (row6) in Job J1, and O22 (row9) in Job J2 were wrongly
#    tagged "P1" - correct them to "S".
# ---------------------------------------------------------------------------
$ws.Range("M4").Value = "S"
$ws.Range("M6").Value = "S"
$ws.Range("M9").Value = "S"

# ---------------------------------------------------------------------------
# 2. Add the new Job J3 (operations O31..O35) in rows 12-16, mirroring the
#    layout/formulas used for J1 and J2 above. Also set number formats for
#    any cells that don't already carry a style from the template rows.
# ---------------------------------------------------------------------------

# Row 12: J3 / O31 / Milling
$ws.Range("A12").Value = "J3"
$ws.Range("B12").Value = "O31"
$ws.Range("C12").Value = "Milling"
$ws.Range("D12").NumberFormat = "0.00"
$ws.Range("D12").Value = 2.8
$ws.Range("E12").NumberFormat = "0.00"
$ws.Range("E12").Value = 0.25
$ws.Range("F12").Value = 800
$ws.Range("G12").Value = 43
$ws.Range("H12").Formula = "=MROUND((D12/(E12*(F12/(PI()*G12))))*60,10)"
$ws.Range("I12").FormulaArray = "=_xlfn.IFS(C12=""Milling"",MROUND(350, 10), C12=""Turning"",MROUND(250, 10), C12=""Drilling"",MROUND(80, 10), C12=""Grinding"",MROUND(52, 10), C12=""Boring"",MROUND(88, 10), C12=""Grooving"",MROUND(230, 10), C12=""Broaching"",MROUND(220, 10), C12=""Honing"",MROUND(90, 10))"
$ws.Range("J12").Formula = "=H12+I12"
$ws.Range("L12").Value = "O33"
$ws.Range("M12").NumberFormat = "0.00"
$ws.Range("M12").Value = "S"
$ws.Range("O12").Value = 0

# Row 13: J3 / O32 / Boring
$ws.Range("A13").Value = "J3"
$ws.Range("B13").Value = "O32"
$ws.Range("C13").Value = "Boring"
$ws.Range("D13").Value = 70
$ws.Range("E13").Value = 0.45
$ws.Range("F13").Value = 3800
$ws.Range("G13").Value = 50
$ws.Range("H13").Formula = "=MROUND((D13/(E13*(F13/(PI()*G13))))*60,10)"
$ws.Range("I13").FormulaArray = "=_xlfn.IFS(C13=""Milling"",MROUND(350, 10), C13=""Turning"",MROUND(250, 10), C13=""Drilling"",MROUND(80, 10), C13=""Grinding"",MROUND(52, 10), C13=""Boring"",MROUND(88, 10), C13=""Grooving"",MROUND(230, 10), C13=""Broaching"",MROUND(220, 10), C13=""Honing"",MROUND(90, 10))"
$ws.Range("J13").Formula = "=H13+I13"
$ws.Range("L13").Value = "O33"
$ws.Range("M13").NumberFormat = "0.00"
$ws.Range("M13").Value = "P2"
$ws.Range("O13").Value = 0

# Row 14: J3 / O33 / Grooving
$ws.Range("A14").Value = "J3"
$ws.Range("B14").Value = "O33"
$ws.Range("C14").Value = "Grooving"
$ws.Range("D14").NumberFormat = "0.00"
$ws.Range("D14").Value = 45
$ws.Range("E14").NumberFormat = "0.00"
$ws.Range("E14").Value = 0.4
$ws.Range("F14").Value = 200
$ws.Range("G14").Value = 5
$ws.Range("H14").Formula = "=MROUND((D14/(E14*(F14/(PI()*G14))))*60,10)"
$ws.Range("I14").FormulaArray = "=_xlfn.IFS(C14=""Milling"",MROUND(350, 10), C14=""Turning"",MROUND(250, 10), C14=""Drilling"",MROUND(80, 10), C14=""Grinding"",MROUND(52, 10), C14=""Boring"",MROUND(88, 10), C14=""Grooving"",MROUND(230, 10), C14=""Broaching"",MROUND(220, 10), C14=""Honing"",MROUND(90, 10))"
$ws.Range("J14").Formula = "=H14+I14"
$ws.Range("K14").Value = "O31,O32"
$ws.Range("L14").Value = "O34,O35"
$ws.Range("M14").NumberFormat = "0.00"
$ws.Range("M14").Value = "S"
$ws.Range("O14").Value = 0

# Row 15: J3 / O34 / Honing
$ws.Range("A15").Value = "J3"
$ws.Range("B15").Value = "O34"
$ws.Range("C15").Value = "Honing"
$ws.Range("D15").NumberFormat = "0.00"
$ws.Range("D15").Value = 0.3
$ws.Range("E15").NumberFormat = "0.00"
$ws.Range("E15").Value = 0.02
$ws.Range("F15").Value = 100
$ws.Range("G15").Value = 8
$ws.Range("H15").Formula = "=MROUND((D15/(E15*(F15/(PI()*G15))))*60,10)"
$ws.Range("I15").FormulaArray = "=_xlfn.IFS(C15=""Milling"",MROUND(350, 10), C15=""Turning"",MROUND(250, 10), C15=""Drilling"",MROUND(80, 10), C15=""Grinding"",MROUND(52, 10), C15=""Boring"",MROUND(88, 10), C15=""Grooving"",MROUND(230, 10), C15=""Broaching"",MROUND(220, 10), C15=""Honing"",MROUND(90, 10))"
$ws.Range("J15").Formula = "=H15+I15"
$ws.Range("K15").Value = "O33"
$ws.Range("M15").NumberFormat = "0.00"
$ws.Range("M15").Value = "S"
$ws.Range("O15").Value = 0

# Row 16: J3 / O35 / Drilling
$ws.Range("A16").Value = "J3"
$ws.Range("B16").Value = "O35"
$ws.Range("C16").Value = "Drilling"
$ws.Range("D16").Value = 45
$ws.Range("E16").Value = 0.4
$ws.Range("F16").Value = 3500
$ws.Range("G16").Value = 52
$ws.Range("H16").Formula = "=MROUND((D16/(E16*(F16/(PI()*G16))))*60,10)"
$ws.Range("I16").FormulaArray = "=_xlfn.IFS(C16=""Milling"",MROUND(350, 10), C16=""Turning"",MROUND(250, 10), C16=""Drilling"",MROUND(80, 10), C16=""Grinding"",MROUND(52, 10), C16=""Boring"",MROUND(88, 10), C16=""Grooving"",MROUND(230, 10), C16=""Broaching"",MROUND(220, 10), C16=""Honing"",MROUND(90, 10))"
$ws.Range("J16").Formula = "=H16+I16"
$ws.Range("K16").Value = "O33"
$ws.Range("M16").NumberFormat = "0.00"
$ws.Range("M16").Value = "P2"
$ws.Range("O16").Value = 0

# ---------------------------------------------------------------------------
# 3. Move the selection to O17, matching where the user ended up after
#    entering the new Scheduling_ERT (O column) data for the new rows.
# ---------------------------------------------------------------------------
$ws.Range("O17").Select()
